$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95 is reset back to full remaining days with a new start date.
$ws.Cells.Item(95, 5).Value = 10
$ws.Cells.Item(95, 6).Value = 20260130

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 95 -or $row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value()
    $cell.Value = $current - 1
}
